$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four "SINGLETON:" rows (originally rows 11, 17, 18, 19) which
# correspond to one-off malware samples that were dropped when the
# train/test split & classifier loading were introduced. Delete from the
# bottom up so earlier row numbers remain valid while deleting.
$ws.Rows.Item(19).EntireRow.Delete()
$ws.Rows.Item(18).EntireRow.Delete()
$ws.Rows.Item(17).EntireRow.Delete()
$ws.Rows.Item(11).EntireRow.Delete()
